# Generate Report for Handback
# Updates the localization-status report after a handback run completed:
#  - flips the zh-cn / de-de status from "Ready for handoff" to
#    "Handed back: in sync with en-US" (Overview + both language sheets)
#  - records the latest target file + handback file/datetime for each
#    language on its own sheet, with a hyperlink back to the source .md
#  - widens a few columns so the new, longer text fits (auto-fit)

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"
$mdFileName = "5d5eae22-f549-4ce7-bf34-414228d83089.md"
$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c25edf8fe20d7e6c7d8c7d087c7f176687fae871/e2e/5d5eae22-f549-4ce7-bf34-414228d83089.md"

# --- 1. Status: Overview sheet (zh-cn / de-de columns) ---
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus

# --- 2. Status on each language sheet's own "Status" column ---
$wsZhCn.Range("C2").Value = $newStatus
$wsDeDe.Range("C2").Value = $newStatus

# --- 3. zh-cn sheet: Latest Target File / Latest Handback File / DateTime ---
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdFileName)
$wsZhCn.Range("J2").Value = "5d5eae22-f549-4ce7-bf34-414228d83089.2faf29f3b6f16e2861c5e7670cc7c6cb6ec7066e.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-09-05 15:13:35"

# --- 4. de-de sheet: Latest Target File / Latest Handback File / DateTime ---
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdFileName)
$wsDeDe.Range("J2").Value = "5d5eae22-f549-4ce7-bf34-414228d83089.2faf29f3b6f16e2861c5e7670cc7c6cb6ec7066e.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-09-05 15:13:43"

# --- 5. Auto-fit the columns whose text just got longer ---
$wsOverview.Columns.Item(5).ColumnWidth = 29.2
$wsOverview.Columns.Item(6).ColumnWidth = 29.2

$wsZhCn.Columns.Item(3).ColumnWidth = 29.2
$wsZhCn.Columns.Item(9).ColumnWidth = 39.15
$wsZhCn.Columns.Item(10).ColumnWidth = 39.15

$wsDeDe.Columns.Item(3).ColumnWidth = 29.2
$wsDeDe.Columns.Item(9).ColumnWidth = 39.15
$wsDeDe.Columns.Item(10).ColumnWidth = 39.15
